$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 5 (bb6f6722 file) now reports "Handed back: in sync with en-US" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F5").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 5 status + handback datetime refreshed, error cleared ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C5").Value = "Handed back: in sync with en-US"
$wsZh.Range("K5").Value = "2016-08-31 07:52:36"
$wsZh.Range("P5").Value = ""

# --- de-de sheet: row 5 status + handback datetime refreshed, error cleared ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C5").Value = "Handed back: in sync with en-US"
$wsDe.Range("K5").Value = "2016-08-31 07:52:53"
$wsDe.Range("P5").Value = ""
